$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 5830
$ws.Range("L3").Value = 6355
$ws.Range("L4").Value = 1570
$ws.Range("L5").Value = 379
$ws.Range("L6").Value = 5218
$ws.Range("L7").Value = 19352

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 379
$ws.Range("L3").Value = 453
$ws.Range("L7").Value = 1275

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L3").Value = 171
$ws.Range("L7").Value = 424

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L6").Value = 252
$ws.Range("L7").Value = 876

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 87
$ws.Range("L7").Value = 272

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L6").Value = 199
$ws.Range("L7").Value = 738

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 374

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L6").Value = 73
$ws.Range("L7").Value = 338

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 164
$ws.Range("L7").Value = 621
$ws.Range("L8").Value = 1275
$ws.Range("L10").Value = 129
$ws.Range("L13").Value = 29
$ws.Range("L14").Value = 97
$ws.Range("L19").Value = 528
$ws.Range("L20").Value = 485
$ws.Range("L22").Value = 59
$ws.Range("L29").Value = 1093
$ws.Range("L31").Value = 192
$ws.Range("L33").Value = 876
$ws.Range("L34").Value = 108
$ws.Range("L36").Value = 244
$ws.Range("L37").Value = 738
$ws.Range("L43").Value = 142
$ws.Range("L46").Value = 46
$ws.Range("L47").Value = 134
$ws.Range("L50").Value = 96
$ws.Range("L51").Value = 246
$ws.Range("L52").Value = 402
$ws.Range("L65").Value = 374
$ws.Range("L67").Value = 667
$ws.Range("L73").Value = 159
$ws.Range("L74").Value = 17
$ws.Range("L76").Value = 295
$ws.Range("L79").Value = 536
$ws.Range("L83").Value = 424
$ws.Range("L85").Value = 958
$ws.Range("L88").Value = 203
$ws.Range("L89").Value = 274
$ws.Range("L90").Value = 201
$ws.Range("L91").Value = 260
$ws.Range("L95").Value = 272
$ws.Range("L96").Value = 221
$ws.Range("L98").Value = 102
$ws.Range("L99").Value = 338
$ws.Range("L101").Value = 19352

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 77
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 192

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L6").Value = 152
$ws.Range("L7").Value = 667

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 323
$ws.Range("L3").Value = 423
$ws.Range("L6").Value = 269
$ws.Range("L7").Value = 1093

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 187
$ws.Range("L3").Value = 164
$ws.Range("L4").Value = 25
$ws.Range("L7").Value = 528

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L3").Value = 58
$ws.Range("L6").Value = 133
$ws.Range("L7").Value = 295

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L2").Value = 39
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("L5").Value = 15
$ws.Range("L6").Value = 29

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 82
$ws.Range("L6").Value = 69

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L2").Value = 68
$ws.Range("L7").Value = 221

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L3").Value = 121
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 172
$ws.Range("L6").Value = 143
$ws.Range("L7").Value = 536

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 151
$ws.Range("L6").Value = 117
$ws.Range("L7").Value = 485

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 86
$ws.Range("L7").Value = 244

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 199
$ws.Range("L5").Value = 18
$ws.Range("L6").Value = 148
$ws.Range("L7").Value = 621

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L2").Value = 37
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 96

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 159

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L3").Value = 82
$ws.Range("L6").Value = 76
$ws.Range("L7").Value = 274

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L2").Value = 66
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 76
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 246

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 47
$ws.Range("L7").Value = 142

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 290
$ws.Range("L3").Value = 393
$ws.Range("L7").Value = 958

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("L4").Value = 7
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L6").Value = 107
$ws.Range("L7").Value = 402

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 17
